# DockingBay.xlsx board-orientation update
#
# The shared-strings table dropped the directional "convr"/"convl" tile
# names and replaced them with new "convbr"/"convbl" ones (belt-oriented
# conveyor tiles). The grid on Sheet1 is re-pointed at the new/shuffled
# shared-string set, and the saved view selection moves from H5 to C5
# (and the stale topLeftCell scroll hint is dropped).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("A2").Value = "floor start 7"
$ws.Range("C2").Value = "floor wall 0"
$ws.Range("E2").Value = "floor wall 0 270"
$ws.Range("H2").Value = "floor wall 0 90"
$ws.Range("J2").Value = "floor wall 0"
$ws.Range("L2").Value = "floor start 8"

# Row 3
$ws.Range("B3").Value = "floor wall 270 start 5"
$ws.Range("C3").Value = "floor wall 270"
$ws.Range("J3").Value = "floor wall 90"
$ws.Range("K3").Value = "floor wall 90 start 6"

# Row 4
$ws.Range("A4").Value = "conv 90 1"
$ws.Range("B4").Value = "conv 90 1"
$ws.Range("C4").Value = "convbr 180 1"
$ws.Range("D4").Value = "floor start 3"
$ws.Range("G4").Value = "floor wall 270"
$ws.Range("I4").Value = "floor start 4"
$ws.Range("J4").Value = "convbl 180 1"

# Row 5
$ws.Range("C5").Value = "convbl 90 1"
$ws.Range("D5").Value = "conv 90 1"
$ws.Range("E5").Value = "conv 90 1"
$ws.Range("F5").Value = "floor start 1"
$ws.Range("G5").Value = "floor wall 270 start 2"
$ws.Range("J5").Value = "convbr 270 1"

# Move the saved selection from H5 to C5 and drop the G1 scroll anchor.
$ws.Range("C5").Select() | Out-Null
